$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the Hass/Palta block (rows 491:492),
# pushing the existing data (previously 491-512) down to 493-514.
$ws.Rows("491:492").Insert()

# Row 491: new "Primera" entry (Hass avocado from Peru, tray of 10 kilos)
$ws.Range("A491").Value = 7
$ws.Range("B491").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C491").Value = "Ñuble"
$ws.Range("D491").Value = 44753
$ws.Range("E491").Value = 16
$ws.Range("F491").Value = "Fruta"
$ws.Range("G491").Value = 100106
$ws.Range("H491").Value = "Oleaginosos"
$ws.Range("I491").Value = 100106002
$ws.Range("J491").Value = "Palta"
$ws.Range("K491").Value = "Hass"
$ws.Range("L491").Value = "Primera"
$ws.Range("M491").Value = 120
$ws.Range("N491").Value = 24000
$ws.Range("O491").Value = 25000
$ws.Range("P491").Value = 24500
$ws.Range("Q491").Value = "$/bandeja 10 kilos"
$ws.Range("R491").Value = "Perú"
$ws.Range("S491").Value = 2450
$ws.Range("T491").Value = 10

# Row 492: new "Segunda" entry (Hass avocado from Peru, tray of 10 kilos)
$ws.Range("A492").Value = 7
$ws.Range("B492").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C492").Value = "Ñuble"
$ws.Range("D492").Value = 44753
$ws.Range("E492").Value = 16
$ws.Range("F492").Value = "Fruta"
$ws.Range("G492").Value = 100106
$ws.Range("H492").Value = "Oleaginosos"
$ws.Range("I492").Value = 100106002
$ws.Range("J492").Value = "Palta"
$ws.Range("K492").Value = "Hass"
$ws.Range("L492").Value = "Segunda"
$ws.Range("M492").Value = 120
$ws.Range("N492").Value = 22000
$ws.Range("O492").Value = 23000
$ws.Range("P492").Value = 22500
$ws.Range("Q492").Value = "$/bandeja 10 kilos"
$ws.Range("R492").Value = "Perú"
$ws.Range("S492").Value = 2250
$ws.Range("T492").Value = 10
